$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Text = "Ablauf eines Projektes mit Scrum`rVerstandnis der Rollen und deren Funktion`rVor- & Nachteile der flexiblen Planung`rGrosserer Planungsaufwand`rBesseres Verstandnis fur Tatigkeit und Probleme der Teammitglieder (Daily)`rProblemlosung im Team`rMVP Konzept`rAuseinandersetzung mit neuem Design-Pattern und dessen konkrete Umsetzung`r Zusammenarbeiten mit Github`rVorteile und Schwierigkeiten von Version Control Systems`rVaadin Framework`rDesignaufwand im Vergleich zu anderen Frameworks riesig`rKaum eine Community vorhanden`rNutzliche Komponenten sind kostenpflichtig (Bsp. Vaadin Designer)`r`r"
$sh.TextFrame.AutoSize = 2
